$wb = $excel.ActiveWorkbook
$wsIncident = $wb.Worksheets.Item("Incident data")
$wsDrivers  = $wb.Worksheets.Item("Drivers data")

# ---------------------------------------------------------------------------
# Fill "Drivers data" sheet with the field-name / SQL-type table (A1:B43)
# ---------------------------------------------------------------------------
$rows = @(
    @("Report Number", "VARCHAR(255)"),
    @("Local Case Number", "VARCHAR(255)"),
    @("Agency Name", "VARCHAR(255)"),
    @("ACRS Report Type", "VARCHAR(255)"),
    @("Crash Date/Time", "DATE"),
    @("Route Type", "VARCHAR(255)"),
    @("Road Name", "VARCHAR(255)"),
    @("Cross-Street Type", "VARCHAR(255)"),
    @("Cross-Street Name", "VARCHAR(255)"),
    @("Off-Road Description", "VARCHAR(255)"),
    @("Municipality", "VARCHAR(255)"),
    @("Related Non-Motorist", "VARCHAR(255)"),
    @("Collision Type", "VARCHAR(255)"),
    @("Weather", "VARCHAR(255)"),
    @("Surface Condition", "VARCHAR(255)"),
    @("Light", "VARCHAR(255)"),
    @("Traffic Control", "VARCHAR(255)"),
    @("Driver Substance Abuse", "VARCHAR(255)"),
    @("Non-Motorist Substance Abuse", "VARCHAR(255)"),
    @("Person ID", "VARCHAR(255)"),
    @("Driver At Fault", "VARCHAR(255)"),
    @("Injury Severity", "VARCHAR(255)"),
    @("Circumstance", "VARCHAR(255)"),
    @("Driver Distracted By", "VARCHAR(255)"),
    @("Drivers License State", "VARCHAR(255)"),
    @("Vehicle ID", "VARCHAR(255)"),
    @("Vehicle Damage Extent", "VARCHAR(255)"),
    @("Vehicle First Impact Location", "VARCHAR(255)"),
    @("Vehicle Second Impact Location", "VARCHAR(255)"),
    @("Vehicle Body Type", "VARCHAR(255)"),
    @("Vehicle Movement", "VARCHAR(255)"),
    @("Vehicle Continuing Dir", "VARCHAR(255)"),
    @("Vehicle Going Dir", "VARCHAR(255)"),
    @("Speed Limit", "INT"),
    @("Driverless Vehicle", "VARCHAR(255)"),
    @("Parked Vehicle", "VARCHAR(255)"),
    @("Vehicle Year", "INT"),
    @("Vehicle Make", "VARCHAR(255)"),
    @("Vehicle Model", "VARCHAR(255)"),
    @("Equipment Problems", "VARCHAR(255)"),
    @("Latitude", "VARCHAR(255)"),
    @("Longitude", "INT"),
    @("Location", "VARCHAR(255)")
)

# Column A is populated in full first, then column B - this mirrors the
# authoring order that produced the shared-string table layout (all of the
# field-name labels first, then "DATE" appended last for column B).
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $wsDrivers.Cells.Item($r, 1).Value = $rows[$i][0]
}
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $wsDrivers.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Match the font/style already used for column B on "Incident data" (Calibri,
# cell style index 1) by copying formats across - this reuses the existing
# style instead of minting a new one.
$wsIncident.Range("B1").Copy() | Out-Null
$wsDrivers.Range("B1:B41").PasteSpecial(-4122) | Out-Null
$wsDrivers.Range("B43").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column A width (bestFit-like)
$wsDrivers.Columns("A:A").ColumnWidth = 24.2

# ---------------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------------
$wsDrivers.Range("A1:B43").Select() | Out-Null
$wsDrivers.Range("B43").Activate() | Out-Null

$wsIncident.Range("B2").Select() | Out-Null

$wsDrivers.Activate() | Out-Null
